$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 34
$ws.Range("B4").Value = 4234
$ws.Range("D5").Value = 4
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 3
$ws.Range("B8").Value = 423
$ws.Range("G8").Value = 324
$ws.Range("C9").Value = 4
$ws.Range("E10").Value = 23
$ws.Range("B12").Value = 32423
$ws.Range("C15").Value = 23

$null = $ws.Range("C11").Select()
